$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) AuditLog sheet: append 4 new audit rows (rows 6-9)
# ---------------------------------------------------------------------
$audit = $wb.Worksheets.Item("AuditLog")

$auditRows = @(
    @{
        row = 6
        A = "AUDIT1745882670934"
        B = "Users"
        C = "2"
        D = "UPDATE"
        E = "1"
        F = "admin"
        G = "2025-04-28T23:24:30.934Z"
        H = '{"before":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-24T00:54:58.039Z","status":"active"},"after":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":["read","write"],"modifiedBy":"admin","modifiedAt":"2025-04-28T23:24:30.933Z","status":"active"}}'
        I = "Updated User 2"
    },
    @{
        row = 7
        A = "AUDIT1745882670934"
        B = "Users"
        C = "2"
        D = "UPDATE"
        E = "1"
        F = "admin"
        G = "2025-04-28T23:24:30.934Z"
        H = '{"before":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-24T00:54:58.039Z","status":"active"},"after":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":["read","write"],"modifiedBy":"admin","modifiedAt":"2025-04-28T23:24:30.933Z","status":"active"}}'
        I = "Updated User 2"
    },
    @{
        row = 8
        A = "AUDIT1745882681205"
        B = "Users"
        C = "3"
        D = "UPDATE"
        E = "1"
        F = "admin"
        G = "2025-04-28T23:24:41.205Z"
        H = '{"before":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-24T00:55:18.315Z","status":"active"},"after":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":["read"],"modifiedBy":"admin","modifiedAt":"2025-04-28T23:24:41.205Z","status":"active"}}'
        I = "Updated User 3"
    },
    @{
        row = 9
        A = "AUDIT1745882681205"
        B = "Users"
        C = "3"
        D = "UPDATE"
        E = "1"
        F = "admin"
        G = "2025-04-28T23:24:41.205Z"
        H = '{"before":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-24T00:55:18.315Z","status":"active"},"after":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":["read"],"modifiedBy":"admin","modifiedAt":"2025-04-28T23:24:41.205Z","status":"active"}}'
        I = "Updated User 3"
    }
)

foreach ($r in $auditRows) {
    $rowNum = $r.row

    $audit.Range("A$rowNum").Value = $r.A
    $audit.Range("B$rowNum").Value = $r.B

    # C and E hold purely-numeric strings ("2", "3", "1") in this data set;
    # force them to remain text (matching the rest of the sheet) instead of
    # being auto-coerced to numbers.
    $cCell = $audit.Range("C$rowNum")
    $cCell.NumberFormat = "@"
    $cCell.Value = $r.C

    $audit.Range("D$rowNum").Value = $r.D

    $eCell = $audit.Range("E$rowNum")
    $eCell.NumberFormat = "@"
    $eCell.Value = $r.E

    $audit.Range("F$rowNum").Value = $r.F
    $audit.Range("G$rowNum").Value = $r.G
    $audit.Range("H$rowNum").Value = $r.H
    $audit.Range("I$rowNum").Value = $r.I
}

# ---------------------------------------------------------------------
# 2) Users sheet: add "permissions" text + refresh "modifiedAt" timestamps
#    for the "user" (row 3) and "viewer" (row 4) accounts.
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

$users.Range("F3").Value = "read,write"
$users.Range("H3").Value = "2025-04-28T23:24:30.933Z"

$users.Range("F4").Value = "read"
$users.Range("H4").Value = "2025-04-28T23:24:41.205Z"
